$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings (e.g. "0.9998")
# are not auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.883.30"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "1.803.28"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "309.97"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4407"
$ws.Range("E7").Value = "  +4.19%  "
$ws.Range("D8").Value = "0.3688"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "0.07416"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("D10").Value = "0.8580"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").Value = "20.74"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.801.59"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "6.625"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "92.97"
$ws.Range("E14").Value = "  +3.05%  "
$ws.Range("D15").Value = "0.07066"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "5.270"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "1.0000"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "0.000008690"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").Value = "0.9992"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "14.81"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").Value = "26.904.89"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "5.160"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "10.82"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "1.973"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "151.52"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").Value = "2.218"
$ws.Range("E26").Value = "  -3.27%  "
$ws.Range("D27").Value = "18.37"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").Value = "5.203"
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("D29").Value = "117.46"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "0.08782"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").Value = "0.7427"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").Value = "1.163"
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("D33").Value = "4.478"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").Value = "2.888"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "0.9990"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "1.092"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "0.01966"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").Value = "0.05207"
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("D39").Value = "0.5245"
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("D40").Value = "7.069"
$ws.Range("E40").Value = "  -3.75%  "
$ws.Range("D41").Value = "2.817"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("D42").Value = "0.1682"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("D43").Value = "8.483"
$ws.Range("E43").Value = "  -1.21%  "
$ws.Range("D44").Value = "0.4999"
$ws.Range("E44").Value = "  +5.78%  "
$ws.Range("D45").Value = "2.036"
$ws.Range("E45").Value = "  +5.80%  "
$ws.Range("D46").Value = "10.36"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").Value = "104.31"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").Value = "0.9989"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("D50").Value = "0.06312"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "0.9212"
$ws.Range("E51").Value = "  +2.13%  "

# Restore the default "Normal" style so cell style indices match the original
# (avoids leaving a stray number-format style reference on the cells).
$ws.Range("D2:D51").Style = "Normal"
